# Apply crypto price/volume updates (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while preserving it as plain text, even when the
# string looks numeric (e.g. "6.160" or "0.9997"), and without leaving the
# cell with a different style than it started with.
function Set-TextValue {
    param($Cell, [string]$Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

$ws.Range("D2").Value = "29.173.36"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "1.842.09"
$ws.Range("E3").Value = "  -0.35%  "
Set-TextValue $ws.Range("D4") "0.9997"
$ws.Range("E4").Value = "  +0.05%  "
Set-TextValue $ws.Range("D5") "241.33"
$ws.Range("E5").Value = "  -1.93%  "
Set-TextValue $ws.Range("D6") "0.6874"
$ws.Range("E6").Value = "  -2.22%  "
$ws.Range("E7").Value = "  +0.06%  "
Set-TextValue $ws.Range("D8") "0.2995"
$ws.Range("E8").Value = "  -2.30%  "
Set-TextValue $ws.Range("D9") "0.07455"
$ws.Range("E9").Value = "  -3.36%  "
Set-TextValue $ws.Range("D10") "23.21"
$ws.Range("E10").Value = "  -1.70%  "
Set-TextValue $ws.Range("D11") "0.07659"
$ws.Range("E11").Value = "  -1.95%  "
$ws.Range("D12").Value = "1.846.85"
$ws.Range("E12").Value = "  -0.11%  "
Set-TextValue $ws.Range("D13") "5.058"
$ws.Range("E13").Value = "  -1.61%  "
Set-TextValue $ws.Range("D14") "0.6824"
$ws.Range("E14").Value = "  -0.73%  "
Set-TextValue $ws.Range("D15") "87.19"
$ws.Range("E15").Value = "  -6.54%  "
Set-TextValue $ws.Range("D16") "6.160"
$ws.Range("E16").Value = "  -6.48%  "
$ws.Range("D17").Value = "29.169.78"
$ws.Range("E17").Value = "  -0.07%  "
Set-TextValue $ws.Range("D18") "0.000008177"
$ws.Range("E18").Value = "  -1.77%  "
$ws.Range("D19").Value = "2.084.08"
$ws.Range("E19").Value = "  -0.49%  "
Set-TextValue $ws.Range("D20") "228.85"
$ws.Range("E20").Value = "  -5.48%  "
Set-TextValue $ws.Range("D21") "12.54"
$ws.Range("E21").Value = "  -1.56%  "
Set-TextValue $ws.Range("D22") "0.9996"
$ws.Range("E22").Value = "  -0.02%  "
Set-TextValue $ws.Range("D23") "7.394"
$ws.Range("E23").Value = "  -1.63%  "
$ws.Range("E24").Value = "  +0.07%  "
Set-TextValue $ws.Range("D25") "0.1448"
$ws.Range("E25").Value = "  -3.93%  "
Set-TextValue $ws.Range("D26") "159.39"
$ws.Range("E26").Value = "  +0.08%  "
Set-TextValue $ws.Range("D27") "8.762"
$ws.Range("E27").Value = "  -1.02%  "
Set-TextValue $ws.Range("D28") "18.08"
$ws.Range("E28").Value = "  -1.27%  "
$ws.Range("E29").Value = "  -1.33%  "
Set-TextValue $ws.Range("D30") "4.279"
$ws.Range("E30").Value = "  +1.13%  "
Set-TextValue $ws.Range("D31") "4.141"
$ws.Range("E31").Value = "  -0.90%  "
Set-TextValue $ws.Range("D32") "1.197"
$ws.Range("E32").Value = "  -0.15%  "
Set-TextValue $ws.Range("D33") "0.05276"
$ws.Range("E33").Value = "  +3.08%  "
Set-TextValue $ws.Range("D34") "0.7611"
$ws.Range("E34").Value = "  -3.43%  "
Set-TextValue $ws.Range("D35") "1.850"
$ws.Range("E35").Value = "  -2.28%  "
$ws.Range("E36").Value = "  -1.19%  "
Set-TextValue $ws.Range("D37") "2.687"
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("D38").Value = "1.304.01"
$ws.Range("E38").Value = "  -1.35%  "
Set-TextValue $ws.Range("D39") "0.01833"
$ws.Range("E39").Value = "  -1.94%  "
Set-TextValue $ws.Range("D40") "2.724"
$ws.Range("E40").Value = "  +0.40%  "
Set-TextValue $ws.Range("D41") "0.9371"
$ws.Range("E41").Value = "  -2.90%  "
Set-TextValue $ws.Range("D42") "5.968"
$ws.Range("E42").Value = "  -1.58%  "
Set-TextValue $ws.Range("D43") "104.97"
$ws.Range("E43").Value = "  -1.84%  "
$ws.Range("D45").Value = "1.985.25"
$ws.Range("E45").Value = "  -0.33%  "
Set-TextValue $ws.Range("D46") "0.5194"
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("E47").Value = "  +0.44%  "
$ws.Range("E48").Value = "  -0.79%  "
Set-TextValue $ws.Range("D49") "9.534"
$ws.Range("E49").Value = "  -1.90%  "
Set-TextValue $ws.Range("D50") "1.771"
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D51") "0.05953"
$ws.Range("E51").Value = "  +0.69%  "
